# Delete 700 pM entry from Shobhan's data (Shobhan2023 / SPR) in the
# VEGF:VEGFR2 Kd dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_VEGFR2")

# The side table in columns G:J (ref/method/value) lists each
# individual measurement; rows 12-14 belonged to Shobhan2023 (SPR),
# with values 700, 770, 270. We need to delete the 700 entry and pull
# the remaining two rows up - but only within columns G:J, since
# columns A:D hold an unrelated per-reference summary table (Rouet
# et al., 2005 on row 12 and Lu et al., 2023 on row 13) that must stay
# exactly where it is. A plain row delete would shift that table too,
# so shift the G:J cells up one at a time instead.
$ws.Range("G12").Value = $ws.Range("G13").Value()
$ws.Range("H12").Value = $ws.Range("H13").Value()
$ws.Range("I12").Value = $ws.Range("I13").Value()
$ws.Range("G13").Value = $ws.Range("G14").Value()
$ws.Range("H13").Value = $ws.Range("H14").Value()
$ws.Range("I13").Value = $ws.Range("I14").Value()
$ws.Range("G14:J14").ClearContents()

# Update the AVERAGE/STDEVA summary in row 11 so it only covers the
# two remaining Shobhan2023 measurements (I12:I13) instead of three.
$ws.Range("C11").Formula = "=AVERAGE(I12:I13)"
$ws.Range("D11").Formula = "=STDEVA(I12:I13)/SQRT(2)"

# Reflect that this was the sheet being edited/selected when saved.
$ws.Activate()
$ws.Range("E14").Select()

# The VEGFA165_NRP1 sheet had two cells (B9, B10) carrying a
# leftover/redundant fill flag on their border style. Clearing the
# (already blank) fill normalizes them back to the same plain border
# style used by equivalent cells elsewhere in the workbook.
$ws3 = $wb.Worksheets.Item("VEGFA165_NRP1")
$ws3.Range("B9").Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternNone
$ws3.Range("B10").Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternNone
